$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 3475568.21
$ws.Range("C7").Value = -21.77566731597328
$ws.Range("D7").Value = 3041
$ws.Range("E7").Value = 3041
$ws.Range("F7").Value = 1142.90306149293
$ws.Range("G7").Value = 21.82520210179235
